$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("K26:L31").NumberFormat = "yyyy\-mm\-dd"
